$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to remain plain text (the source data stores
# prices/percentages as literal strings, not numbers), then write the
# updated values scraped on the new run date.
$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","E25","D26","E26","E27","E28","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D47","E47","D48","E48")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "258.36"
$ws.Range("E2").Value = "1.34%"
$ws.Range("D3").Value = "26.89"
$ws.Range("E3").Value = "-4.14%"
$ws.Range("D4").Value = "4.776"
$ws.Range("E4").Value = "-10.81%"
$ws.Range("D5").Value = "0.05969"
$ws.Range("E5").Value = "2.51%"
$ws.Range("D6").Value = "6.684"
$ws.Range("E6").Value = "-0.41%"
$ws.Range("D7").Value = "0.8772"
$ws.Range("E7").Value = "1.65%"
$ws.Range("D8").Value = "0.9540"
$ws.Range("E8").Value = "4.58%"
$ws.Range("D9").Value = "0.1417"
$ws.Range("E9").Value = "-0.34%"
$ws.Range("D10").Value = "0.03610"
$ws.Range("E10").Value = "4.37%"
$ws.Range("D11").Value = "0.07207"
$ws.Range("E11").Value = "0.57%"
$ws.Range("D12").Value = "0.03143"
$ws.Range("E12").Value = "-1.97%"
$ws.Range("D13").Value = "0.09238"
$ws.Range("E13").Value = "-0.02%"
$ws.Range("D14").Value = "0.001541"
$ws.Range("E14").Value = "0.07%"
$ws.Range("D15").Value = "0.0006078"
$ws.Range("E15").Value = "0.30%"
$ws.Range("D16").Value = "0.005966"
$ws.Range("E16").Value = "1.37%"
$ws.Range("D17").Value = "3.486"
$ws.Range("E17").Value = "-0.29%"
$ws.Range("D18").Value = "3.227"
$ws.Range("E18").Value = "0.02%"
$ws.Range("D19").Value = "2.218"
$ws.Range("E19").Value = "-1.51%"
$ws.Range("D20").Value = "0.3135"
$ws.Range("E20").Value = "-1.01%"
$ws.Range("E21").Value = "-1.07%"
$ws.Range("D22").Value = "3.534"
$ws.Range("E22").Value = "0.18%"
$ws.Range("D23").Value = "0.04216"
$ws.Range("E23").Value = "1.54%"
$ws.Range("E25").Value = "-0.14%"
$ws.Range("D26").Value = "0.004515"
$ws.Range("E26").Value = "-11.81%"
$ws.Range("E27").Value = "-0.02%"
$ws.Range("E28").Value = "-22.99%"
$ws.Range("D40").Value = "0.03848"
$ws.Range("E40").Value = "0.20%"
$ws.Range("D41").Value = "0.005988"
$ws.Range("E41").Value = "4.89%"
$ws.Range("D42").Value = "0.1105"
$ws.Range("E42").Value = "0.60%"
$ws.Range("D43").Value = "0.002299"
$ws.Range("E43").Value = "4.53%"
$ws.Range("D44").Value = "0.01108"
$ws.Range("E44").Value = "10.84%"
$ws.Range("D45").Value = "0.00005491"
$ws.Range("E45").Value = "3.81%"
$ws.Range("D47").Value = "0.08548"
$ws.Range("E47").Value = "-14.53%"
$ws.Range("D48").Value = "0.002123"
$ws.Range("E48").Value = "-4.07%"
